# Updated naming and js for multiple entries per row.
# Adds a header row with "Total Number of Units" (col B) and "Notes" (col C)
# above the existing two label rows, and moves the active selection to F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Total Number of Units"
$ws.Range("C1").Value = "Notes"

$ws.Range("F9").Select()
